$d = $word.ActiveDocument

# --- 1) Simplify the M3 paragraph: merge the split/proofErr-wrapped runs into a single run ---
$d.Content.Find.Execute(
    "M3: Each exchange must be defined by a unique identifier of the form ex:p:q where p and q are numbers between 0 and 9",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "M3: Each exchange must be defined by a unique identifier of the form ex:p:q where p and q are numbers between 0 and 9",
    2) | Out-Null

# --- 2) Insert a brand-new "M6: The user must be able to input the id and the
#         location of the exchange" paragraph right after M5, before the old M6 paragraph ---
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "M5: Customer must be at*") {
        $p.Range.InsertParagraphAfter() | Out-Null
        break
    }
}

foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -eq "`r") {
        $prev = $p.Previous()
        $next = $p.Next()
        $isTarget = ($prev -ne $null) -and ($next -ne $null) `
            -and ($prev.Range.Text -like "M5: Customer must be at*") `
            -and ($next.Range.Text -like "M6: User must be able to input the current capacity*")
        if ($isTarget) {
            $p.Range.Text = "M6: The user must be able to input the id and the location of the exchange"
            break
        }
    }
}

# --- 3) Turn the old "M6: User must be able to input the current capacity of a
#         given exchange" paragraph into "M7: ..." splitting "M7" into its own run,
#         with the (single, special) _GoBack bookmark relocated right after it ---
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "M6: User must be able to input the current capacity of a given exchange*") {
        $pStart = $p.Range.Start

        $prefixRange = $d.Range($pStart, $pStart + 2)
        $prefixRange.Text = "M7"

        $bmPos = $pStart + 2
        $bmRange = $d.Range($bmPos, $bmPos)
        $d.Bookmarks.Add("_GoBack", $bmRange) | Out-Null
        break
    }
}
